# Regenerate orders with updated distance/size codes.
# The experiment's distance conditions and one size condition were renamed:
#   D80 -> D86
#   D64 -> D69
#   D51 -> D55
#   S30 -> S31
# These codes appear embedded inside many string values across the sheet
# (Condition, Filename_Left, Filename_Right, Distance, Size columns), so we
# walk every used cell and rewrite any text value containing the old codes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
  for ($c = 1; $c -le $colCount; $c++) {
    $cell = $ws.Cells.Item($r, $c)
    $val = $cell.Value()
    if ($val -is [string]) {
      $newVal = $val -replace 'D80','D86' -replace 'D64','D69' -replace 'D51','D55' -replace 'S30','S31'
      if ($newVal -ne $val) {
        $cell.Value = $newVal
      }
    }
  }
}
